$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, pushing existing rows 28-39 down to 29-40.
$ws.Rows.Item(28).Insert()

# Copy the number format (date style) used in column D from the row below (now row 29)
# into the new row 28's D cell so the date renders consistently.
$ws.Range("D29").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 28 with the inserted record's values.
$ws.Cells.Item(28, 1).Value = 4
$ws.Cells.Item(28, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(28, 3).Value = "Los Lagos"
$ws.Cells.Item(28, 4).Value = 44463
$ws.Cells.Item(28, 5).Value = 10
$ws.Cells.Item(28, 6).Value = 100112026
$ws.Cells.Item(28, 7).Value = "Haba"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 80
$ws.Cells.Item(28, 11).Value = 16000
$ws.Cells.Item(28, 12).Value = 16000
$ws.Cells.Item(28, 13).Value = 16000
$ws.Cells.Item(28, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(28, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(28, 16).Value = 640
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"
